# Applies:
#  1) "ind_process_routes_fom" (sheet10) and "ind_process_routes_capex" (sheet9):
#     - For data rows 2-78, set the 2030/2040/2050 columns (G,H,I) equal to the
#       2025 column (F) value -> costs are now expressed flat in EUR-2025.
#     - For rows 79-80 (refineries (H2)FT-DC / (H2)MeOH-DC) all four year
#       columns (F,G,H,I) are replaced with freshly recalculated EUR-2025 values.
#  2) "ind_production_2018_nuts1" (sheet3) and "ind_production_30_50_nuts1"
#     (sheet4): remove the rows whose production value(s) are zero.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Re-base cost sheets onto the 2025 (EUR-2025) value for all future years
# ---------------------------------------------------------------------------
$costSheetNames = @("ind_process_routes_fom", "ind_process_routes_capex")

foreach ($sheetName in $costSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($r = 2; $r -le 78; $r++) {
        $baseValue = $ws.Cells.Item($r, 6).Value2
        $ws.Cells.Item($r, 7).Value2 = $baseValue
        $ws.Cells.Item($r, 8).Value2 = $baseValue
        $ws.Cells.Item($r, 9).Value2 = $baseValue
    }
}

# Row 79 / 80 on each sheet get brand new (recomputed) values for every year
$fom = $wb.Worksheets.Item("ind_process_routes_fom")
$fom.Cells.Item(79, 6).Value2 = 19.28455846613214
$fom.Cells.Item(79, 7).Value2 = 16.54537904625254
$fom.Cells.Item(79, 8).Value2 = 11.07927604282393
$fom.Cells.Item(79, 9).Value2 = 9.645343192148706

$fom.Cells.Item(80, 6).Value2 = 5.456365489631739
$fom.Cells.Item(80, 7).Value2 = 4.197204222793643
$fom.Cells.Item(80, 8).Value2 = 4.197204222793643
$fom.Cells.Item(80, 9).Value2 = 3.637576993087825

$capex = $wb.Worksheets.Item("ind_process_routes_capex")
$capex.Cells.Item(79, 6).Value2 = 275.1267368041234
$capex.Cells.Item(79, 7).Value2 = 237.8415726249732
$capex.Cells.Item(79, 8).Value2 = 163.6909646889521
$capex.Cells.Item(79, 9).Value2 = 134.3105351293966

$capex.Cells.Item(80, 6).Value2 = 188.874190025714
$capex.Cells.Item(80, 7).Value2 = 152.4984200948358
$capex.Cells.Item(80, 8).Value2 = 134.3105351293966
$capex.Cells.Item(80, 9).Value2 = 121.7189224610157

# ---------------------------------------------------------------------------
# 2) Drop rows with zero production from the NUTS1 production sheets
# ---------------------------------------------------------------------------
# Rows (1-based, matching the original sheet layout) that have zero
# production in every value column. Deleted from the bottom up so row
# numbers of rows still to be removed stay valid.
$zeroProductionRows = @(208, 204, 202, 200, 197, 196, 193, 192, 191, 189, 188)

$prod2018 = $wb.Worksheets.Item("ind_production_2018_nuts1")
foreach ($r in $zeroProductionRows) {
    $prod2018.Rows.Item($r).Delete()
}

$prod3050 = $wb.Worksheets.Item("ind_production_30_50_nuts1")
foreach ($r in $zeroProductionRows) {
    $prod3050.Rows.Item($r).Delete()
}
